$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Antworten")

$ws.Range("I11").Formula = "=AVERAGE(I2:I9)"
$ws.Range("I12").Formula = "=STDEVA(I2:I9)"

$ws.Columns.Item(9).ColumnWidth = 7.90625
